# Auto-generated edit script updating cryptos list values (commit: "Updated cryptos list on Sun Dec 31 19:43:36 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.918.29'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '2.293.63'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'316.20"
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').Value = "'104.14"
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').Value = "'0.624"
$ws.Range('E7').Value = '  -1.11%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('D10').Value = "'39.38"
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').Value = "'8.46"
$ws.Range('E12').Value = '  +1.44%  '
$ws.Range('E13').Value = '  +2.09%  '
$ws.Range('E14').Value = '  +4.30%  '
$ws.Range('D15').Value = "'15.32"
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').Value = '2.640.10'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').Value = '2.303.28'
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('D18').Value = '42.813.83'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').Value = "'7.45"
$ws.Range('E19').Value = '  -0.93%  '
$ws.Range('D20').Value = "'13.88"
$ws.Range('E20').Value = '  +25.62%  '
$ws.Range('E21').Value = '  -0.77%  '
$ws.Range('D22').Value = "'73.98"
$ws.Range('E22').Value = '  +0.93%  '
$ws.Range('D23').Value = "'3.55"
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').Value = "'263.25"
$ws.Range('D25').Value = "'2.20"
$ws.Range('E25').Value = '  -3.57%  '
$ws.Range('E26').Value = '  +0.41%  '
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('E28').Value = '  -1.83%  '
$ws.Range('D29').Value = "'7.03"
$ws.Range('E29').Value = '  +19.13%  '
$ws.Range('E30').Value = '  -1.83%  '
$ws.Range('D31').Value = "'37.44"
$ws.Range('E31').Value = '  +4.80%  '
$ws.Range('D32').Value = "'166.49"
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('D33').Value = "'0.0872"
$ws.Range('E33').Value = '  -0.38%  '
$ws.Range('E34').Value = '  -4.65%  '
$ws.Range('E35').Value = '  -0.64%  '
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('E37').Value = '  -0.55%  '
$ws.Range('E38').Value = '  -5.82%  '
$ws.Range('D39').Value = "'3.83"
$ws.Range('E39').Value = '  +2.60%  '
$ws.Range('E40').Value = '  -2.83%  '
$ws.Range('D41').Value = "'1.58"
$ws.Range('E41').Value = '  +4.56%  '
$ws.Range('E42').Value = '  +0.98%  '
$ws.Range('D43').Value = "'69.52"
$ws.Range('E43').Value = '  -0.96%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').Value = "'92.49"
$ws.Range('E45').Value = '  -1.95%  '
$ws.Range('D46').Value = "'12.15"
$ws.Range('E46').Value = '  +0.91%  '
$ws.Range('D47').Value = "'113.95"
$ws.Range('E47').Value = '  +0.73%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.725.25'
$ws.Range('E48').Value = '  +8.36%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').Value = "'80.21"
$ws.Range('E49').Value = '  -3.02%  '
$ws.Range('E50').Value = '  -1.61%  '
$ws.Range('D51').Value = "'5.15"
$ws.Range('E51').Value = '  +1.35%  '
